# Fix the pump dose volume for row 2 (Pump 1 / F.0.1.22_1): the dose
# volume in column G ("Dose vol.") was set too low, which caused the
# connection/handshake with the pump to fail. Correct it to 300.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 300

# Move/leave the active selection where the author left it when saving.
$ws.Range("K9").Select()
